$wb = $excel.ActiveWorkbook

# The handoff attempt failed before a target file was produced, so the
# per-language "Ready for handoff" status becomes "Handoff transform failed"
# everywhere that text is used (Overview rollup + each language sheet).
$ws0 = $wb.Worksheets.Item("Overview")
$ws0.Range("B2").Value = "Handoff transform failed"
$ws0.Range("C2").Value = "Handoff transform failed"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status (B2): "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # No handoff file was produced, so "Latest Handoff File" (C2), which
    # held a hyperlinked filename, is cleared out entirely.
    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime (D2) resets to the zero/default datetime,
    # matching the other unset datetime cells (D3/G2/G3).
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason (H2): "Include" -> "Ignored"
    $ws.Range("H2").Value = "Ignored"
}
